## Fix Training Data Issue (#48)
## The BF column ("Date") held a malformed/backwards date string
## ("6-8-2011-12") for every data row. NBA stats for the game date of
## June 8, 2012 were shown one day off, so the date needs to be corrected
## to the ISO form "2012-06-08" for every data row (BF2:BF31 -- BF1 is
## just the "Date" header and is left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow  = 31
$dateColumn   = "BF"
$correctDate  = "2012-06-08"

# Write the corrected date as a text-producing formula first. Assigning
# "2012-06-08" straight to .Value would make Excel "smart" parse the
# yyyy-mm-dd look-alike into a real date serial, which is not what the
# source data should contain -- it needs to stay the literal text string.
for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {
    $ws.Range($dateColumn + $row).Formula = '="' + $correctDate + '"'
}

# Flatten the formulas down to plain literal values (copy / paste-special
# values-only) so the cells end up holding the text "2012-06-08" itself,
# with no left-over formula and no change to the cells' existing
# (unstyled) formatting.
$dataRange = $ws.Range($dateColumn + $firstDataRow + ":" + $dateColumn + $lastDataRow)
$dataRange.Copy()
$dataRange.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
